$d = $word.ActiveDocument

# --- Paragraph 1 (heading): merge "First Edition, " + "First Correction"
#     into a single run "First Edition, Second Correction" ---
$p1 = $d.Paragraphs(2)
$p1.Range.Find.Execute("First Edition, First Correction", $true, $false, $false, `
    $false, $false, $true, 1, $false, "First Edition, Second Correction", 2)

# --- Paragraph 2 (body): turn the single run containing
#     "...Terms, First Edition, First Correction (...1e1c)." into five runs:
#       "...Terms, " | "First" | " Edition, " | "Second" | " Correction (...1e2c)."
$p2 = $d.Paragraphs(3)
$full = $p2.Range
$startPos = $full.Start

$seg1 = "We agree to these terms, along with the Square One Standard Contractor Terms, "
$seg2 = "First"
$seg3 = " Edition, "
$seg4 = "First"
$seg5 = " Correction (https://squareoneforms.com/contractor/1e1c)."

$off1 = 0
$off2 = $off1 + $seg1.Length
$off3 = $off2 + $seg2.Length
$off4 = $off3 + $seg3.Length
$off5 = $off4 + $seg4.Length
$off6 = $off5 + $seg5.Length

# Step 1: apply the text edits first (while the paragraph is still one run),
# using offsets computed from the known original text.

# 4th segment: "First" -> "Second"
$r4 = $d.Range($startPos + $off4, $startPos + $off5)
$r4.Text = "Second"
$delta = ("Second").Length - $seg4.Length
$off5b = $off5 + $delta
$off6b = $off6 + $delta

# 5th segment: update the edition code in the URL from 1e1c to 1e2c
$r5 = $d.Range($startPos + $off5b, $startPos + $off6b)
$r5.Find.Execute("1e1c", $true, $false, $false, $false, $false, $true, 1, $false, "1e2c", 2)

# Step 2: split the single run into five runs at the segment boundaries by
# toggling Bold on and back off across each inner span -- this forces Word
# to create separate <w:r> elements without altering the final formatting.
$rA = $d.Range($startPos + $off2, $startPos + $off3)
$rA.Bold = 1
$rA.Bold = 0

$rB = $d.Range($startPos + $off3, $startPos + $off4)
$rB.Bold = 1
$rB.Bold = 0

$rC = $d.Range($startPos + $off4, $startPos + $off5b)
$rC.Bold = 1
$rC.Bold = 0
